$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reword D1's label: spell out the eligible professor ranks instead of the
# generic "enter a number" instruction. E1 ("國籍") is untouched.
$ws.Range("D1").Value = "外籍學者身分（教授、副教授、助理教授或博士後研究員）"

# Column D needs to be much wider now that it holds the longer label
# (target stored width is 55.625 characters; feed the COM layer a value
# that quantizes back to the closest representable width).
$ws.Columns.Item(4).ColumnWidth = 54.86

# Reset the view: zoom 160% -> 115%, and move the active selection from
# B12 to E4.
$excel.ActiveWindow.Zoom = 115
$ws.Range("E4").Select()
